$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 68; this shifts the existing rows 68-183 down
# to 69-184 (and the sheet's used range grows to A1:R184).
$ws.Rows(68).Insert()

# Populate the newly inserted row 68 with this week's data point.
$ws.Cells.Item(68, 1).Value = 7
$ws.Cells.Item(68, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(68, 3).Value = "Ñuble"
$ws.Cells.Item(68, 4).Value = 44477
$ws.Cells.Item(68, 5).Value = 16
$ws.Cells.Item(68, 6).Value = 100114013
$ws.Cells.Item(68, 7).Value = "Zanahoria"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 120
$ws.Cells.Item(68, 11).Value = 8000
$ws.Cells.Item(68, 12).Value = 9000
$ws.Cells.Item(68, 13).Value = 8500
$ws.Cells.Item(68, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(68, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(68, 16).Value = 425
$ws.Cells.Item(68, 17).Value = 20
$ws.Cells.Item(68, 18).Value = "Hortaliza"
